$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.259.54'
$ws.Range('E2').Value = '  -2.04%  '
$ws.Range('D3').Value = '3.596.53'
$ws.Range('E3').Value = '  -1.13%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'582.50"
$ws.Range('E5').Value = '  -1.35%  '
$ws.Range('D6').Value = "'173.78"
$ws.Range('E6').Value = '  -3.56%  '
$ws.Range('D7').Value = "'0.628"
$ws.Range('E7').Value = '  +2.63%  '
$ws.Range('D8').Value = '3.588.99'
$ws.Range('E8').Value = '  -1.06%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  -4.50%  '
$ws.Range('D11').Value = "'6.57"
$ws.Range('E11').Value = '  +12.63%  '
$ws.Range('D12').Value = "'0.611"
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').Value = "'48.04"
$ws.Range('E13').Value = '  -3.39%  '
$ws.Range('E14').Value = '  -1.88%  '
$ws.Range('D15').Value = "'692.13"
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('D16').Value = '4.183.82'
$ws.Range('E16').Value = '  -0.89%  '
$ws.Range('D17').Value = "'8.96"
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').Value = '70.411.68'
$ws.Range('E18').Value = '  -2.03%  '
$ws.Range('D19').Value = '3.608.48'
$ws.Range('E19').Value = '  -1.95%  '
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('D21').Value = "'17.61"
$ws.Range('E21').Value = '  -3.86%  '
$ws.Range('D22').Value = "'11.35"
$ws.Range('E22').Value = '  -2.22%  '
$ws.Range('D23').Value = "'0.927"
$ws.Range('E23').Value = '  -0.84%  '
$ws.Range('D24').Value = "'16.93"
$ws.Range('E24').Value = '  -4.75%  '
$ws.Range('D25').Value = "'98.94"
$ws.Range('E25').Value = '  -4.28%  '
$ws.Range('D26').Value = "'3.89"
$ws.Range('E26').Value = '  -3.29%  '
$ws.Range('D27').Value = "'2.73"
$ws.Range('E27').Value = '  -3.78%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = "'9.60"
$ws.Range('E29').Value = '  -3.87%  '
$ws.Range('D30').Value = "'34.19"
$ws.Range('E30').Value = '  -2.31%  '
$ws.Range('D31').Value = "'9.08"
$ws.Range('E31').Value = '  -1.16%  '
$ws.Range('E32').Value = '  -5.30%  '
$ws.Range('D33').Value = "'7.44"
$ws.Range('E33').Value = '  +2.58%  '
$ws.Range('D34').Value = "'1.36"
$ws.Range('E34').Value = '  -4.97%  '
$ws.Range('D35').Value = "'3.90"
$ws.Range('E35').Value = '  -6.98%  '
$ws.Range('D36').Value = "'575.92"
$ws.Range('E36').Value = '  -0.75%  '
$ws.Range('D37').Value = "'11.02"
$ws.Range('E37').Value = '  -2.78%  '
$ws.Range('E38').Value = '  -2.73%  '
$ws.Range('D39').Value = "'58.20"
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('D40').Value = "'1.00"
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('E41').Value = '  -4.29%  '
$ws.Range('E42').Value = '  -1.65%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '3.515.28'
$ws.Range('E43').Value = '  -4.36%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').Value = "'0.343"
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('D45').Value = "'33.92"
$ws.Range('E45').Value = '  -5.09%  '
$ws.Range('D46').Value = '0.0₃0720'
$ws.Range('E46').Value = '  -5.55%  '
$ws.Range('D47').Value = "'2.95"
$ws.Range('E47').Value = '  +4.99%  '
$ws.Range('D48').Value = "'2.61"
$ws.Range('E48').Value = '  -5.12%  '
$ws.Range('E49').Value = '  +1.65%  '
$ws.Range('D50').Value = "'135.47"
$ws.Range('E50').Value = '  +2.71%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.151"
$ws.Range('E51').Value = '  +1.25%  '
